# Automatische test-sync: 2025-08-30 19:26:50
#
# Adds the newly logged "Opvolging contact" mail entry to the Logs sheet
# (row 8) and refreshes the aggregated counts / ordering on the Dashboard
# sheet to reflect the new total for "Klantenservice / Contact".

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new mail-log row ------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A8").Value = "Opvolging contact"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("D8").Value = "Klantenservice / Contact"
$logs.Range("F8").Value = "2025-08-30 19:26:33"
$logs.Range("G8").Value = "Nee"
$logs.Range("H8").Value = "Ja"
$logs.Range("I8").Value = "Nee"
$logs.Range("J8").Value = "Nee"

# Extend the conditional-formatting ranges (D/G/H/I/J) so the new row 8 is
# covered too, just like the other data rows. All rules that share a sqref
# get updated together when the first rule's range is modified.
$logs.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D8")) | Out-Null
$logs.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G8")) | Out-Null
$logs.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H8")) | Out-Null
$logs.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I8")) | Out-Null
$logs.Range("J2:J7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J8")) | Out-Null

# --- Dashboard sheet: refresh the per-category summary --------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

# "Klantenservice / Contact" now has 2 mails (was 1) and moves up to row 4,
# while "Planning / Afspraak" (still 1 mail) drops to row 5.
$dashboard.Range("A4").Value = "Klantenservice / Contact"
$dashboard.Range("B4").Value = 2
$dashboard.Range("A5").Value = "Planning / Afspraak"
$dashboard.Range("B5").Value = 1
